$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 10 new rows above the old row 2 (nov25), pushing nov25/dec25/old-row4
# down to rows 12/13/14. This preserves the existing shared-string references
# for nov25 (row -> 12) and dec25 (row -> 13) untouched.
$ws.Rows("2:11").Insert()

# --- Header row: A1 ("month") switches from Arial/General to Arial/Text format.
$ws.Range("A1").NumberFormat = "@"

# --- New row 2: jan25 (same look as the old nov25/dec25 rows: Arial 10, text/number fmt)
$ws.Range("A2").Value = "jan25"
$ws.Range("A2").Font.Name = "Arial"
$ws.Range("A2").Font.Size = 10
$ws.Range("A2").NumberFormat = "@"
$ws.Range("B2").Value = 996240.32
$ws.Range("B2").Font.Name = "Arial"
$ws.Range("B2").Font.Size = 10
$ws.Range("B2").NumberFormat = "#,##0.00"
$ws.Range("C2").Value = 376.82
$ws.Range("C2").Font.Name = "Arial"
$ws.Range("C2").Font.Size = 10
$ws.Range("D2").Value = 2114856
$ws.Range("D2").Font.Name = "Arial"
$ws.Range("D2").Font.Size = 10
$ws.Range("D2").NumberFormat = "#,##0.00"

# --- New row 3: feb25 (same formatting family as row 2)
$ws.Range("A3").Value = "feb25"
$ws.Range("A3").Font.Name = "Arial"
$ws.Range("A3").Font.Size = 10
$ws.Range("A3").NumberFormat = "@"
$ws.Range("B3").Value = 890205.59
$ws.Range("B3").Font.Name = "Arial"
$ws.Range("B3").Font.Size = 10
$ws.Range("B3").NumberFormat = "#,##0.00"
$ws.Range("C3").Value = 505.93
$ws.Range("C3").Font.Name = "Arial"
$ws.Range("C3").Font.Size = 10
$ws.Range("D3").Value = 2114856
$ws.Range("D3").Font.Name = "Arial"
$ws.Range("D3").Font.Size = 10
$ws.Range("D3").NumberFormat = "#,##0.00"

# --- New row 4: mar25 (default/Aptos font from here on, D keeps Arial/number fmt)
$ws.Range("A4").Value = "mar25"
$ws.Range("A4").NumberFormat = "@"
$ws.Range("B4").Value = 910973.82
$ws.Range("B4").NumberFormat = "#,##0.00"
$ws.Range("C4").Value = 554.42999999999995
$ws.Range("D4").Value = 2114856
$ws.Range("D4").Font.Name = "Arial"
$ws.Range("D4").Font.Size = 10
$ws.Range("D4").NumberFormat = "#,##0.00"

# --- New row 5: apr25
$ws.Range("A5").Value = "apr25"
$ws.Range("A5").NumberFormat = "@"
$ws.Range("B5").Value = 935486.9
$ws.Range("B5").NumberFormat = "#,##0.00"
$ws.Range("C5").Value = 946.24
$ws.Range("C5").Font.Name = "Arial"
$ws.Range("C5").Font.Size = 10
$ws.Range("D5").Value = 2114856
$ws.Range("D5").Font.Name = "Arial"
$ws.Range("D5").Font.Size = 10
$ws.Range("D5").NumberFormat = "#,##0.00"

# --- New row 6: may25
$ws.Range("A6").Value = "may25"
$ws.Range("A6").NumberFormat = "@"
$ws.Range("B6").Value = 940898.39
$ws.Range("B6").NumberFormat = "#,##0.00"
$ws.Range("C6").Value = 1114.1099999999999
$ws.Range("C6").Font.Name = "Arial"
$ws.Range("C6").Font.Size = 10
$ws.Range("D6").Value = 2114856
$ws.Range("D6").Font.Name = "Arial"
$ws.Range("D6").Font.Size = 10
$ws.Range("D6").NumberFormat = "#,##0.00"

# --- New row 7: jun25
$ws.Range("A7").Value = "jun25"
$ws.Range("A7").NumberFormat = "@"
$ws.Range("B7").Value = 911165.81
$ws.Range("B7").NumberFormat = "#,##0.00"
$ws.Range("C7").Value = 937.85
$ws.Range("C7").Font.Name = "Arial"
$ws.Range("C7").Font.Size = 10
$ws.Range("D7").Value = 2114856
$ws.Range("D7").Font.Name = "Arial"
$ws.Range("D7").Font.Size = 10
$ws.Range("D7").NumberFormat = "#,##0.00"

# --- New row 8: jul25 (D switches to the new 2487916.6 total + default/Aptos font)
$ws.Range("A8").Value = "jul25"
$ws.Range("A8").NumberFormat = "@"
$ws.Range("B8").Value = 914121.36
$ws.Range("B8").NumberFormat = "#,##0.00"
$ws.Range("C8").Value = 814.39
$ws.Range("D8").Value = 2487916.6
$ws.Range("D8").NumberFormat = "#,##0.00"

# --- New row 9: aug25
$ws.Range("A9").Value = "aug25"
$ws.Range("A9").NumberFormat = "@"
$ws.Range("B9").Value = 976007.07
$ws.Range("B9").NumberFormat = "#,##0.00"
$ws.Range("C9").Value = 847.5
$ws.Range("D9").Value = 2487916.6
$ws.Range("D9").NumberFormat = "#,##0.00"

# --- New row 10: sep25
$ws.Range("A10").Value = "sep25"
$ws.Range("A10").NumberFormat = "@"
$ws.Range("B10").Value = 984294.02
$ws.Range("B10").NumberFormat = "#,##0.00"
$ws.Range("C10").Value = 1002.42
$ws.Range("C10").NumberFormat = "#,##0.00"
$ws.Range("D10").Value = 2487916.6
$ws.Range("D10").NumberFormat = "#,##0.00"

# --- New row 11: oct25
$ws.Range("A11").Value = "oct25"
$ws.Range("A11").NumberFormat = "@"
$ws.Range("B11").Value = 1064340.31
$ws.Range("B11").NumberFormat = "#,##0.00"
$ws.Range("C11").Value = 763.79
$ws.Range("D11").Value = 2487916.6
$ws.Range("D11").NumberFormat = "#,##0.00"

# --- Remove the old row 4 (it shifted to row 14 and carried stray thick-bottom
# border / D4 formatting that doesn't belong in the new layout).
$ws.Rows.Item(14).Delete()

# --- The old row 3 (dec25) shifted to row 13 and still carries its original
# "thick bottom border" row formatting; clear that residual row formatting.
$ws.Rows.Item(13).AutoFit()

# --- Two new, empty but formatted rows at the bottom (21/22), matching the
# look of the other Arial-formatted rows in the sheet.
$ws.Range("A21:D21").Font.Name = "Arial"
$ws.Range("A21:D21").Font.Size = 10
$ws.Range("A21").NumberFormat = "@"

$ws.Range("A22:D22").Font.Name = "Arial"
$ws.Range("A22:D22").Font.Size = 10
$ws.Range("A22").NumberFormat = "@"
$ws.Range("B22").NumberFormat = "#,##0.00"

# --- Restore the selection/active cell shown in the saved view.
$ws.Range("D19").Select()
